# Auto-generated edit script applying the Omega_Profits market-data refresh
$wb = $excel.ActiveWorkbook

# --- ALC (Worksheets.Item(1)) ---
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(17, 8).Value = 3169.9412
$ws.Cells.Item(17, 10).Value = 3169.9412
$ws.Cells.Item(17, 12).Value = 9509.8236
$ws.Cells.Item(17, 14).Value = -9845.8236
$ws.Cells.Item(70, 8).Value = 168899.5
$ws.Cells.Item(70, 10).Value = 334133
$ws.Cells.Item(70, 12).Value = 1002399
$ws.Cells.Item(70, 14).Value = -1002939
$ws.Cells.Item(73, 8).Value = 168899.5
$ws.Cells.Item(73, 10).Value = 334133
$ws.Cells.Item(73, 12).Value = 1002399
$ws.Cells.Item(73, 14).Value = -1004271
$ws.Cells.Item(118, 8).Value = 1303.6
$ws.Cells.Item(118, 9).Value = 1303.6
$ws.Cells.Item(118, 11).Value = 3910.8
$ws.Cells.Item(118, 13).Value = -2253.8
$ws.Cells.Item(121, 8).Value = 1404.4
$ws.Cells.Item(121, 10).Value = 1404.4
$ws.Cells.Item(121, 12).Value = 4213.200000000001
$ws.Cells.Item(121, 14).Value = -7707.200000000001
$ws.Cells.Item(132, 8).Value = 3464.2368
$ws.Cells.Item(132, 9).Value = 3075.457
$ws.Cells.Item(132, 11).Value = 9226.370999999999
$ws.Cells.Item(132, 13).Value = -6696.370999999999
$ws.Cells.Item(137, 8).Value = 2194.8667
$ws.Cells.Item(137, 9).Value = 2156.2307
$ws.Cells.Item(137, 11).Value = 6468.6921
$ws.Cells.Item(137, 13).Value = -3918.6921
$ws.Cells.Item(138, 8).Value = 2857.4329
$ws.Cells.Item(138, 9).Value = 1725.1765
$ws.Cells.Item(138, 10).Value = 3242.4
$ws.Cells.Item(138, 11).Value = 5175.529500000001
$ws.Cells.Item(138, 12).Value = 9727.200000000001
$ws.Cells.Item(138, 13).Value = -35.52950000000055
$ws.Cells.Item(138, 14).Value = -20007.2

# --- ARM (Worksheets.Item(2)) ---
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(21, 8).Value = 9757.5
$ws.Cells.Item(21, 9).Value = 5015
$ws.Cells.Item(21, 10).Value = 14500
$ws.Cells.Item(21, 11).Value = 5015
$ws.Cells.Item(21, 12).Value = 14500
$ws.Cells.Item(21, 14).Value = -15248
$ws.Cells.Item(32, 8).Value = 7095.0225
$ws.Cells.Item(32, 9).Value = 2979.3142
$ws.Cells.Item(32, 11).Value = 2979.3142
$ws.Cells.Item(32, 13).Value = -2692.3142
$ws.Cells.Item(63, 8).Value = 6117.1665
$ws.Cells.Item(63, 9).Value = 2269.875
$ws.Cells.Item(63, 10).Value = 9195
$ws.Cells.Item(63, 11).Value = 2269.875
$ws.Cells.Item(63, 12).Value = 9195
$ws.Cells.Item(63, 13).Value = -1583.875
$ws.Cells.Item(63, 14).Value = -10567
$ws.Cells.Item(66, 8).Value = 6117.1665
$ws.Cells.Item(66, 9).Value = 2269.875
$ws.Cells.Item(66, 10).Value = 9195
$ws.Cells.Item(66, 11).Value = 11349.375
$ws.Cells.Item(66, 12).Value = 45975
$ws.Cells.Item(66, 13).Value = -7917.375
$ws.Cells.Item(66, 14).Value = -52839
$ws.Cells.Item(74, 8).Value = 2177.5
$ws.Cells.Item(74, 9).Value = 2197.3684
$ws.Cells.Item(74, 11).Value = 2197.3684
$ws.Cells.Item(74, 13).Value = -1323.3684
$ws.Cells.Item(77, 8).Value = 2177.5
$ws.Cells.Item(77, 9).Value = 2197.3684
$ws.Cells.Item(77, 11).Value = 10986.842
$ws.Cells.Item(77, 13).Value = -6618.841999999999
$ws.Cells.Item(132, 8).Value = 1804.1082
$ws.Cells.Item(132, 9).Value = 1553.4062
$ws.Cells.Item(132, 11).Value = 4660.2186
$ws.Cells.Item(132, 13).Value = -2130.2186

# --- BSM (Worksheets.Item(3)) ---
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(80, 8).Value = 705.0625
$ws.Cells.Item(80, 10).Value = 837.6923
$ws.Cells.Item(80, 12).Value = 837.6923
$ws.Cells.Item(80, 14).Value = -2833.6923
$ws.Cells.Item(83, 8).Value = 705.0625
$ws.Cells.Item(83, 10).Value = 837.6923
$ws.Cells.Item(83, 12).Value = 4188.4615
$ws.Cells.Item(83, 14).Value = -14172.4615
$ws.Cells.Item(86, 8).Value = 2550.2856
$ws.Cells.Item(86, 9).Value = 2290.5
$ws.Cells.Item(86, 10).Value = 3199.75
$ws.Cells.Item(86, 11).Value = 2290.5
$ws.Cells.Item(86, 12).Value = 3199.75
$ws.Cells.Item(86, 13).Value = -1167.5
$ws.Cells.Item(86, 14).Value = -5445.75
$ws.Cells.Item(89, 8).Value = 2550.2856
$ws.Cells.Item(89, 9).Value = 2290.5
$ws.Cells.Item(89, 10).Value = 3199.75
$ws.Cells.Item(89, 11).Value = 11452.5
$ws.Cells.Item(89, 12).Value = 15998.75
$ws.Cells.Item(89, 13).Value = -5836.5
$ws.Cells.Item(89, 14).Value = -27230.75
$ws.Cells.Item(107, 8).Value = 2111.7273
$ws.Cells.Item(107, 9).Value = 2022.9
$ws.Cells.Item(107, 11).Value = 2022.9
$ws.Cells.Item(107, 13).Value = -102.9000000000001

# --- CRP (Worksheets.Item(4)) ---
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(31, 8).Value = 6472.9165
$ws.Cells.Item(31, 9).Value = 9629.210999999999
$ws.Cells.Item(31, 11).Value = 9629.210999999999
$ws.Cells.Item(31, 13).Value = -9334.210999999999
$ws.Cells.Item(34, 8).Value = 6472.9165
$ws.Cells.Item(34, 9).Value = 9629.210999999999
$ws.Cells.Item(34, 11).Value = 9629.210999999999
$ws.Cells.Item(34, 13).Value = -9427.210999999999
$ws.Cells.Item(58, 8).Value = 2187.4285
$ws.Cells.Item(58, 9).Value = 2099.3572
$ws.Cells.Item(58, 11).Value = 2099.3572
$ws.Cells.Item(58, 13).Value = -1896.3572
$ws.Cells.Item(64, 8).Value = 105000
$ws.Cells.Item(64, 10).Value = 105000
$ws.Cells.Item(64, 12).Value = 105000
$ws.Cells.Item(64, 14).Value = -105496
$ws.Cells.Item(67, 8).Value = 105000
$ws.Cells.Item(67, 10).Value = 105000
$ws.Cells.Item(67, 12).Value = 105000
$ws.Cells.Item(67, 14).Value = -106716
$ws.Cells.Item(86, 8).Value = 25661536
$ws.Cells.Item(86, 9).Value = 25661536
$ws.Cells.Item(86, 11).Value = 25661536
$ws.Cells.Item(86, 13).Value = -25660413
$ws.Cells.Item(87, 8).Value = 20000
$ws.Cells.Item(87, 9).Value = 20000
$ws.Cells.Item(87, 11).Value = 20000
$ws.Cells.Item(87, 13).Value = -18814
$ws.Cells.Item(89, 8).Value = 25661536
$ws.Cells.Item(89, 9).Value = 25661536
$ws.Cells.Item(89, 11).Value = 128307680
$ws.Cells.Item(89, 13).Value = -128302064
$ws.Cells.Item(90, 8).Value = 20000
$ws.Cells.Item(90, 9).Value = 20000
$ws.Cells.Item(90, 11).Value = 60000
$ws.Cells.Item(90, 13).Value = -54072
$ws.Cells.Item(99, 8).Value = 6838.294
$ws.Cells.Item(99, 9).Value = 4328.6
$ws.Cells.Item(99, 11).Value = 4328.6
$ws.Cells.Item(99, 13).Value = -2830.6
$ws.Cells.Item(107, 8).Value = 31251404
$ws.Cells.Item(107, 9).Value = 41668220
$ws.Cells.Item(107, 11).Value = 41668220
$ws.Cells.Item(107, 13).Value = -41666300
$ws.Cells.Item(126, 8).Value = 6838.294
$ws.Cells.Item(126, 9).Value = 4328.6
$ws.Cells.Item(126, 11).Value = 12985.8
$ws.Cells.Item(126, 13).Value = -10515.8
$ws.Cells.Item(132, 8).Value = 1998
$ws.Cells.Item(132, 9).Value = 1958
$ws.Cells.Item(132, 11).Value = 5874
$ws.Cells.Item(132, 13).Value = -3344
$ws.Cells.Item(134, 8).Value = 1263.5625
$ws.Cells.Item(134, 9).Value = 1202.0714
$ws.Cells.Item(134, 11).Value = 3606.2142
$ws.Cells.Item(134, 13).Value = -1071.2142
$ws.Cells.Item(136, 8).Value = 2187.4285
$ws.Cells.Item(136, 9).Value = 2099.3572
$ws.Cells.Item(136, 11).Value = 6298.071599999999
$ws.Cells.Item(136, 13).Value = -3748.071599999999
$ws.Cells.Item(141, 8).Value = 294224.88
$ws.Cells.Item(141, 10).Value = 294224.88
$ws.Cells.Item(141, 12).Value = 294224.88
$ws.Cells.Item(141, 14).Value = -304584.88

# --- CUL (Worksheets.Item(5)) ---
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(3, 8).Value = 4209.875
$ws.Cells.Item(3, 9).Value = 4209.875
$ws.Cells.Item(3, 11).Value = 12629.625
$ws.Cells.Item(3, 13).Value = -12517.625
$ws.Cells.Item(4, 8).Value = 88474530
$ws.Cells.Item(4, 9).Value = 333563360
$ws.Cells.Item(4, 10).Value = 6778244.5
$ws.Cells.Item(4, 11).Value = 1000690080
$ws.Cells.Item(4, 12).Value = 20334733.5
$ws.Cells.Item(4, 13).Value = -1000689968
$ws.Cells.Item(4, 14).Value = -20334957.5
$ws.Cells.Item(5, 8).Value = 478.1579
$ws.Cells.Item(5, 9).Value = 403.33334
$ws.Cells.Item(5, 11).Value = 1210.00002
$ws.Cells.Item(5, 13).Value = -1098.00002
$ws.Cells.Item(122, 8).Value = 2272
$ws.Cells.Item(122, 10).Value = 3253.1428
$ws.Cells.Item(122, 12).Value = 29278.2852
$ws.Cells.Item(122, 14).Value = -34178.2852
$ws.Cells.Item(131, 8).Value = 2176.52
$ws.Cells.Item(131, 10).Value = 2542.2354
$ws.Cells.Item(131, 12).Value = 7626.706200000001
$ws.Cells.Item(131, 14).Value = -17706.7062
$ws.Cells.Item(135, 8).Value = 478.1579
$ws.Cells.Item(135, 9).Value = 403.33334
$ws.Cells.Item(135, 11).Value = 3630.00006
$ws.Cells.Item(135, 13).Value = -1095.00006
$ws.Cells.Item(139, 8).Value = 3182.7222
$ws.Cells.Item(139, 9).Value = 3018.0625
$ws.Cells.Item(139, 10).Value = 4500
$ws.Cells.Item(139, 11).Value = 9054.1875
$ws.Cells.Item(139, 12).Value = 13500
$ws.Cells.Item(139, 13).Value = -3914.1875
$ws.Cells.Item(139, 14).Value = -23780
$ws.Cells.Item(141, 8).Value = 5839
$ws.Cells.Item(141, 9).Value = 5839
$ws.Cells.Item(141, 10).Value = 0
$ws.Cells.Item(141, 11).Value = 17517
$ws.Cells.Item(141, 12).Value = 0
$ws.Cells.Item(141, 13).Value = -12337

# --- GSM (Worksheets.Item(6)) ---
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(99, 8).Value = 20391.143
$ws.Cells.Item(99, 9).Value = 12123.167
$ws.Cells.Item(99, 11).Value = 12123.167
$ws.Cells.Item(99, 13).Value = -9877.166999999999
$ws.Cells.Item(132, 8).Value = 5204.52
$ws.Cells.Item(132, 9).Value = 4929.2607
$ws.Cells.Item(132, 10).Value = 8370
$ws.Cells.Item(132, 11).Value = 14787.7821
$ws.Cells.Item(132, 12).Value = 25110
$ws.Cells.Item(132, 13).Value = -12257.7821
$ws.Cells.Item(132, 14).Value = -30170
$ws.Cells.Item(136, 8).Value = 31456.385
$ws.Cells.Item(136, 10).Value = 31456.385
$ws.Cells.Item(136, 12).Value = 94369.155
$ws.Cells.Item(136, 14).Value = -99469.155

# --- LTW (Worksheets.Item(7)) ---
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(40, 8).Value = 5489.8945
$ws.Cells.Item(40, 9).Value = 5324.8125
$ws.Cells.Item(40, 11).Value = 5324.8125
$ws.Cells.Item(40, 13).Value = -5188.8125
$ws.Cells.Item(82, 8).Value = 1394.125
$ws.Cells.Item(82, 9).Value = 1142.5454
$ws.Cells.Item(82, 10).Value = 1947.6
$ws.Cells.Item(82, 11).Value = 1142.5454
$ws.Cells.Item(82, 12).Value = 1947.6
$ws.Cells.Item(82, 13).Value = -781.5454
$ws.Cells.Item(82, 14).Value = -2669.6
$ws.Cells.Item(85, 8).Value = 1394.125
$ws.Cells.Item(85, 9).Value = 1142.5454
$ws.Cells.Item(85, 10).Value = 1947.6
$ws.Cells.Item(85, 11).Value = 1142.5454
$ws.Cells.Item(85, 12).Value = 1947.6
$ws.Cells.Item(85, 13).Value = 105.4546
$ws.Cells.Item(85, 14).Value = -4443.6
$ws.Cells.Item(122, 8).Value = 11473
$ws.Cells.Item(122, 9).Value = 11473
$ws.Cells.Item(122, 11).Value = 34419
$ws.Cells.Item(122, 13).Value = -31969

# Row 141 on CUL loses its N141 cell entirely (was -22360, now removed)
$wsCUL = $wb.Worksheets.Item(5)
$wsCUL.Cells.Item(141, 14).ClearContents()
